$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume/percentage change (E) columns
# Values are forced to text via a leading quote prefix (mirrors how the
# source data is stored as text in the workbook) and the style is reset
# to "Normal" afterwards so no stray quote-prefix formatting is left behind.

$ws.Range("D2").Value = "'26.764.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.33%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.642.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.08%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.39%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'218.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.48%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.22%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.253"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.31%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.40%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.08%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.869.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.19%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.632.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.07%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E15").Value = "'  -0.33%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'64.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.45%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'26.770.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.27%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.0₃0735"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.05%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'215.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.13%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.20%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.95%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +6.35%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.29%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -2.08%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'145.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.30%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.14%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.62%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.40%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.32%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.96%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +1.38%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +1.04%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.44%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.286.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.39%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.19%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +1.45%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -0.39%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +1.50%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.817"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.18%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +0.26%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.61%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -1.62%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'5.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.47%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.780.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.14%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'61.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.52%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'91.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.30%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.47%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.28%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'7.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.63%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.17%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.03%  "
$ws.Range("E51").Style = "Normal"
